$wb = $excel.ActiveWorkbook

# ALC row 28 (Leve Item ID 27772)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 243.125
$ws.Range("I28").Value = 243.125
$ws.Range("K28").Value = 243.125
$ws.Range("M28").Value = 241.875

# ALC row 98 (Leve Item ID 36237)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1256.5454
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# ALC row 113 (Leve Item ID 27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 5946.1113
$ws.Range("I113").Value = 4402.8184
$ws.Range("K113").Value = 4402.8184
$ws.Range("M113").Value = -1148.8184

# ALC row 122 (Leve Item ID 36237)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1256.5454
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# ARM row 45 (Leve Item ID 27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 863
$ws.Range("I45").Value = 900
$ws.Range("K45").Value = 900
$ws.Range("M45").Value = -523

# ARM row 110 (Leve Item ID 27708)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3049
$ws.Range("J110").Value = 4999.6665
$ws.Range("L110").Value = 4999.6665
$ws.Range("N110").Value = -9089.666499999999

# BSM row 48 (Leve Item ID 22888)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

# BSM row 99 (Leve Item ID 19943)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4600.2
$ws.Range("I99").Value = 4600.2
$ws.Range("K99").Value = 4600.2
$ws.Range("M99").Value = -3102.2

# CRP row 36 (Leve Item ID 1845)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

# CRP row 40 (Leve Item ID 1845)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

# CRP row 70 (Leve Item ID 12011)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 21000
$ws.Range("J70").Value = 21000
$ws.Range("L70").Value = 21000
$ws.Range("N70").Value = -21630

# CRP row 73 (Leve Item ID 12011)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 21000
$ws.Range("J73").Value = 21000
$ws.Range("L73").Value = 21000
$ws.Range("N73").Value = -23184

# CRP row 74 (Leve Item ID 10636)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 29999.857
$ws.Range("J74").Value = 29999.857
$ws.Range("L74").Value = 29999.857
$ws.Range("N74").Value = -31747.857

# CRP row 77 (Leve Item ID 10636)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 29999.857
$ws.Range("J77").Value = 29999.857
$ws.Range("L77").Value = 89999.571
$ws.Range("N77").Value = -98735.571

# CRP row 105 (Leve Item ID 19928)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1084.2858
$ws.Range("I105").Value = 1015
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 1015
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 732
$ws.Range("N105").Value = -4994

# CRP row 134 (Leve Item ID 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1306
$ws.Range("I134").Value = 1306
$ws.Range("K134").Value = 3918
$ws.Range("M134").Value = -1383

# CUL row 4 (Leve Item ID 4650)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2178.238
$ws.Range("I4").Value = 2204
$ws.Range("J4").Value = 2143.889
$ws.Range("K4").Value = 6612
$ws.Range("L4").Value = 6431.667
$ws.Range("M4").Value = -6500
$ws.Range("N4").Value = -6655.667

# CUL row 33 (Leve Item ID 4867)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 259
$ws.Range("I33").Value = 174.75
$ws.Range("J33").Value = 596
$ws.Range("K33").Value = 1048.5
$ws.Range("L33").Value = 3576
$ws.Range("M33").Value = -765.5
$ws.Range("N33").Value = -4142

# CUL row 36 (Leve Item ID 4732)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 200
$ws.Range("I36").Value = 200
$ws.Range("K36").Value = 600
$ws.Range("M36").Value = -431

# CUL row 41 (Leve Item ID 4700)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 2000
$ws.Range("J41").Value = 2000
$ws.Range("L41").Value = 6000
$ws.Range("N41").Value = -6676

# CUL row 141 (Leve Item ID 44076)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 2249.5
$ws.Range("I141").Value = 2249.5
$ws.Range("K141").Value = 6748.5
$ws.Range("M141").Value = -1568.5

# GSM row 70 (Leve Item ID 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 333333340
$ws.Range("I70").Value = 333333340
$ws.Range("K70").Value = 333333340
$ws.Range("M70").Value = -333333070

# GSM row 73 (Leve Item ID 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 333333340
$ws.Range("I73").Value = 333333340
$ws.Range("K73").Value = 333333340
$ws.Range("M73").Value = -333332404

# GSM row 99 (Leve Item ID 19532)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 16985.5
$ws.Range("I99").Value = 14314
$ws.Range("J99").Value = 25000
$ws.Range("K99").Value = 14314
$ws.Range("L99").Value = 25000
$ws.Range("M99").Value = -12068
$ws.Range("N99").Value = -29492

# GSM row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1046.2222
$ws.Range("I122").Value = 1061.2
$ws.Range("J122").Value = 1027.5
$ws.Range("K122").Value = 3183.6
$ws.Range("L122").Value = 3082.5
$ws.Range("M122").Value = -733.6000000000004
$ws.Range("N122").Value = -7982.5

# GSM row 126 (Leve Item ID 36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 15582
$ws.Range("I126").Value = 12637
$ws.Range("J126").Value = 19999.5
$ws.Range("K126").Value = 37911
$ws.Range("L126").Value = 59998.5
$ws.Range("M126").Value = -35441
$ws.Range("N126").Value = -64938.5

# LTW row 16 (Leve Item ID 5289)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 622
$ws.Range("I16").Value = 622
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 622
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -452
$ws.Range("N16").ClearContents()

# LTW row 40 (Leve Item ID 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5746.5835
$ws.Range("I40").Value = 5359.909
$ws.Range("K40").Value = 5359.909
$ws.Range("M40").Value = -5223.909

# LTW row 68 (Leve Item ID 12563)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2998.3333
$ws.Range("I68").Value = 2750
$ws.Range("K68").Value = 2750
$ws.Range("M68").Value = -2001

# LTW row 71 (Leve Item ID 12563)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2998.3333
$ws.Range("I71").Value = 2750
$ws.Range("K71").Value = 13750
$ws.Range("M71").Value = -10006

# LTW row 106 (Leve Item ID 18713)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

# LTW row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1499.5
$ws.Range("I136").Value = 1499.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4498.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1948.5
$ws.Range("N136").ClearContents()

# WVR row 70 (Leve Item ID 11979)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 32221.666
$ws.Range("J70").Value = 32221.666
$ws.Range("L70").Value = 32221.666
$ws.Range("N70").Value = -32851.666

# WVR row 73 (Leve Item ID 11979)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 32221.666
$ws.Range("J73").Value = 32221.666
$ws.Range("L73").Value = 32221.666
$ws.Range("N73").Value = -34405.666

# WVR row 75 (Leve Item ID 11957)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

# WVR row 78 (Leve Item ID 11957)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

# WVR row 122 (Leve Item ID 36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2386.3333
$ws.Range("I122").Value = 2386.3333
$ws.Range("K122").Value = 7158.999899999999
$ws.Range("M122").Value = -4708.999899999999
